# Generate Report for Handoff
# Refreshes the localization-status report: updates the "Latest HO Xliff
# Generate Date" for the batch of four e2e\*.md files (rows 4-7) on the
# Overview sheet, and updates the Priority ("low" -> "ht") and Latest
# Handoff Datetime for those same four rows on both the zh-cn and de-de
# localization sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
foreach ($r in 4..7) {
    $overview.Cells.Item($r, 7).Value = "2016-09-03 04:35:09"
}

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-03 04:34:59"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($r in 4..7) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-09-03 04:35:09"
}
